$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ17472117",
    "summ17732574",
    "summ18075329",
    "summ18430403",
    "summ18766068",
    "summ19138498",
    "summ19417185",
    "summ19740887",
    "summ20054334",
    "summ20367269",
    "summ20660523",
    "summ20962757",
    "summ21261819",
    "summ21541695",
    "summ21850942",
    "summ22147901",
    "summ22465431",
    "summ22764240",
    "summ23067778",
    "summ23357748",
    "summ23655230",
    "summ23998542",
    "summ24298843",
    "summ24593920",
    "summ24897452",
    "summ25213047",
    "summ25519641",
    "summ25810947",
    "summ26104540",
    "summ26395617",
    "summ26700961",
    "summ27001802",
    "summ27278046",
    "summ27583061",
    "summ27870460",
    "summ28161991",
    "summ28471537",
    "summ28775998",
    "summ29072060",
    "summ29390812",
    "summ29733601",
    "summ30080537",
    "summ30444501",
    "summ30790328",
    "summ31107817",
    "summ31446039",
    "summ31757442",
    "summ32055339",
    "summ32357711",
    "summ32664054"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}

